$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 2022.5454
$ws.Range("I17").Value2 = 1845
$ws.Range("J17").Value2 = 2062
$ws.Range("K17").Value2 = 5535
$ws.Range("L17").Value2 = 6186
$ws.Range("M17").Value2 = -5367
$ws.Range("N17").Value2 = -6522
$ws.Range("H97").Value2 = 2354.125
$ws.Range("J97").Value2 = 2354.125
$ws.Range("L97").Value2 = 7062.375
$ws.Range("N97").Value2 = -8054.375
$ws.Range("H112").Value2 = 2239.4167
$ws.Range("J112").Value2 = 2239.4167
$ws.Range("L112").Value2 = 6718.250100000001
$ws.Range("N112").Value2 = -8934.250100000001
$ws.Range("H116").Value2 = 6774.5225
$ws.Range("I116").Value2 = 6824.4585
$ws.Range("J116").Value2 = 6746.6514
$ws.Range("K116").Value2 = 6824.4585
$ws.Range("L116").Value2 = 6746.6514
$ws.Range("M116").Value2 = -3382.4585
$ws.Range("N116").Value2 = -13630.6514
$ws.Range("H121").Value2 = 4150
$ws.Range("J121").Value2 = 4150
$ws.Range("L121").Value2 = 12450
$ws.Range("N121").Value2 = -15944
$ws.Range("H132").Value2 = 6650.3784
$ws.Range("I132").Value2 = 6945.9395
$ws.Range("J132").Value2 = 4212
$ws.Range("K132").Value2 = 20837.8185
$ws.Range("L132").Value2 = 12636
$ws.Range("M132").Value2 = -18307.8185
$ws.Range("N132").Value2 = -17696
$ws.Range("H137").Value2 = 5004083.5
$ws.Range("I137").Value2 = 10003278
$ws.Range("J137").Value2 = 4889
$ws.Range("K137").Value2 = 30009834
$ws.Range("L137").Value2 = 14667
$ws.Range("M137").Value2 = -30007284
$ws.Range("N137").Value2 = -19767
$ws.Range("H138").Value2 = 3238.7646
$ws.Range("I138").Value2 = 3286
$ws.Range("J138").Value2 = 3219.0833
$ws.Range("K138").Value2 = 9858
$ws.Range("L138").Value2 = 9657.249899999999
$ws.Range("M138").Value2 = -4718
$ws.Range("N138").Value2 = -19937.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2812.3088
$ws.Range("I32").Value2 = 2552.2742
$ws.Range("J32").Value2 = 5499.3335
$ws.Range("K32").Value2 = 2552.2742
$ws.Range("L32").Value2 = 5499.3335
$ws.Range("M32").Value2 = -2265.2742
$ws.Range("N32").Value2 = -6073.3335
$ws.Range("H61").Value2 = 2472.9565
$ws.Range("I61").Value2 = 2179.8948
$ws.Range("J61").Value2 = 3865
$ws.Range("K61").Value2 = 2179.8948
$ws.Range("L61").Value2 = 3865
$ws.Range("M61").Value2 = -1967.8948
$ws.Range("N61").Value2 = -4289
$ws.Range("H132").Value2 = 4446702.5
$ws.Range("I132").Value2 = 1774.7413
$ws.Range("J132").Value2 = 19611750
$ws.Range("K132").Value2 = 5324.2239
$ws.Range("L132").Value2 = 58835250
$ws.Range("M132").Value2 = -2794.2239
$ws.Range("N132").Value2 = -58840310
$ws.Range("H136").Value2 = 2472.9565
$ws.Range("I136").Value2 = 2179.8948
$ws.Range("J136").Value2 = 3865
$ws.Range("K136").Value2 = 6539.6844
$ws.Range("L136").Value2 = 11595
$ws.Range("M136").Value2 = -3989.6844
$ws.Range("N136").Value2 = -16695

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1149.0869
$ws.Range("I107").Value2 = 1105.3636
$ws.Range("J107").Value2 = 2111
$ws.Range("K107").Value2 = 1105.3636
$ws.Range("L107").Value2 = 2111
$ws.Range("M107").Value2 = 814.6364000000001
$ws.Range("N107").Value2 = -5951
$ws.Range("H134").Value2 = 2485.0557
$ws.Range("I134").Value2 = 2199.6086
$ws.Range("J134").Value2 = 4126.375
$ws.Range("K134").Value2 = 6598.825800000001
$ws.Range("L134").Value2 = 12379.125
$ws.Range("M134").Value2 = -4063.825800000001
$ws.Range("N134").Value2 = -17449.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1238.909
$ws.Range("I22").Value2 = 1506.3334
$ws.Range("J22").Value2 = 918
$ws.Range("K22").Value2 = 1506.3334
$ws.Range("L22").Value2 = 918
$ws.Range("M22").Value2 = -1156.3334
$ws.Range("N22").Value2 = -1618
$ws.Range("H31").Value2 = 5346.077
$ws.Range("I31").Value2 = 5067.3335
$ws.Range("K31").Value2 = 5067.3335
$ws.Range("M31").Value2 = -4772.3335
$ws.Range("H34").Value2 = 5346.077
$ws.Range("I34").Value2 = 5067.3335
$ws.Range("K34").Value2 = 5067.3335
$ws.Range("M34").Value2 = -4865.3335
$ws.Range("H58").Value2 = 1863.2727
$ws.Range("I58").Value2 = 1306.1333
$ws.Range("J58").Value2 = 3057.1428
$ws.Range("K58").Value2 = 1306.1333
$ws.Range("L58").Value2 = 3057.1428
$ws.Range("M58").Value2 = -1103.1333
$ws.Range("N58").Value2 = -3463.1428
$ws.Range("H94").Value2 = 2042
$ws.Range("I94").Value2 = 1912.25
$ws.Range("J94").Value2 = 2171.75
$ws.Range("K94").Value2 = 1912.25
$ws.Range("L94").Value2 = 2171.75
$ws.Range("M94").Value2 = -1461.25
$ws.Range("N94").Value2 = -3073.75
$ws.Range("H134").Value2 = 2693.6758
$ws.Range("I134").Value2 = 2728.5334
$ws.Range("J134").Value2 = 2544.2856
$ws.Range("K134").Value2 = 8185.600199999999
$ws.Range("L134").Value2 = 7632.8568
$ws.Range("M134").Value2 = -5650.600199999999
$ws.Range("N134").Value2 = -12702.8568
$ws.Range("H136").Value2 = 1863.2727
$ws.Range("I136").Value2 = 1306.1333
$ws.Range("J136").Value2 = 3057.1428
$ws.Range("K136").Value2 = 3918.3999
$ws.Range("L136").Value2 = 9171.428400000001
$ws.Range("M136").Value2 = -1368.3999
$ws.Range("N136").Value2 = -14271.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 1391.7894
$ws.Range("I122").Value2 = 534.3333
$ws.Range("J122").Value2 = 1787.5385
$ws.Range("K122").Value2 = 4808.9997
$ws.Range("L122").Value2 = 16087.8465
$ws.Range("M122").Value2 = -2358.9997
$ws.Range("N122").Value2 = -20987.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 3050.5715
$ws.Range("I132").Value2 = 2073.4546
$ws.Range("J132").Value2 = 6633.3335
$ws.Range("K132").Value2 = 6220.3638
$ws.Range("L132").Value2 = 19900.0005
$ws.Range("M132").Value2 = -3690.3638
$ws.Range("N132").Value2 = -24960.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value2 = 302
$ws.Range("J115").Value2 = 302
$ws.Range("L115").Value2 = 302
$ws.Range("N115").Value2 = -2652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 10104093
$ws.Range("I132").Value2 = 14495686
$ws.Range("J132").Value2 = 3429.7
$ws.Range("K132").Value2 = 43487058
$ws.Range("L132").Value2 = 10289.1
$ws.Range("M132").Value2 = -43484528
$ws.Range("N132").Value2 = -15349.1
$ws.Range("H136").Value2 = 6274.4683
$ws.Range("I136").Value2 = 6180.641
$ws.Range("J136").Value2 = 6731.875
$ws.Range("K136").Value2 = 18541.923
$ws.Range("L136").Value2 = 20195.625
$ws.Range("M136").Value2 = -15991.923
$ws.Range("N136").Value2 = -25295.625
